# 3_1_Journal_Barras_Matias.xlsx — add the missing "Sprint 3" journal entry
# and move the viewport/selection onto the row that was just filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal_Barras_Matias")

# Row 22 is the first data row of the (until now empty) third sprint block.
# Fill in the task description and the hours spent on it.
$ws.Range("B22").Value = "Création du scripte pour extraire le texte et teste d'autres solution (écheque)"
$ws.Range("D22").Value = 7

# The sprint-3 subtotal (D27 = SUM(D22:D26)) and the grand total (D85) are
# formulas already on the sheet, so they will pick up the new value once
# Excel recalculates.
$excel.Calculate()

# Move the active selection onto the row that was edited, matching where
# the author was working when they saved the file.
$ws.Activate()
$ws.Range("B22:C22").Select()
